$wb = $excel.ActiveWorkbook

$msfl = $wb.Worksheets.Item("msfl")
$wsfl = $wb.Worksheets.Item("wsfl")

# Apply an AutoFilter on the msfl sheet (A1:C77), filtering the
# "Position" column (column C, the 3rd column of the filter range)
# down to just "Defender". Passing the criteria as an array produces a
# discrete values filter (<filters><filter val="Defender"/></filters>)
# rather than a custom-operator filter.
$msfl.Range("A1:C77").AutoFilter(3, @("Defender"))

# Re-create the wsfl!_FilterDatabase defined name after adding msfl's,
# so the workbook-level definedNames collection ends up ordered by
# sheet index (msfl=0 before wsfl=1), matching a fresh AutoFilter save.
$wb.Names.Item("wsfl!_FilterDatabase").Delete()

$msfl.Names.Add("_xlnm._FilterDatabase", "=msfl!`$A`$1:`$C`$77") | Out-Null
$wsfl.Names.Add("_xlnm._FilterDatabase", "=wsfl!`$C`$1:`$C`$34") | Out-Null

# Builtin _FilterDatabase names are hidden; mark both as such (re-fetch
# by name so the Visible assignment actually sticks).
$wb.Names.Item("msfl!_FilterDatabase").Visible = $false
$wb.Names.Item("wsfl!_FilterDatabase").Visible = $false

# Make msfl the active sheet/tab (this also clears wsfl's tabSelected
# flag since only one sheet can be the active tab) and set the
# selection to the visible filtered data range.
$msfl.Activate()
$msfl.Range("A8:C77").Select()
